$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue 2 4 "72.679.49"
Set-TextValue 2 5 "  +5.70%  "
Set-TextValue 3 4 "4.058.88"
Set-TextValue 3 5 "  +5.50%  "
Set-TextValue 4 5 "  +0.13%  "
Set-TextValue 5 4 "522.51"
Set-TextValue 5 5 "  -0.14%  "
Set-TextValue 6 4 "148.15"
Set-TextValue 6 5 "  +4.10%  "
Set-TextValue 7 4 "0.728"
Set-TextValue 7 5 "  +20.24%  "
Set-TextValue 8 4 "4.051.17"
Set-TextValue 8 5 "  +5.57%  "
Set-TextValue 9 5 "  +0.16%  "
Set-TextValue 10 4 "0.785"
Set-TextValue 10 5 "  +10.39%  "
Set-TextValue 11 5 "  +4.67%  "
Set-TextValue 12 4 "0.0000333"
Set-TextValue 12 5 "  +0.88%  "
Set-TextValue 13 4 "48.36"
Set-TextValue 13 5 "  +16.28%  "
Set-TextValue 14 4 "11.15"
Set-TextValue 14 5 "  +10.12%  "
Set-TextValue 15 4 "4.698.44"
Set-TextValue 15 5 "  +5.33%  "
Set-TextValue 16 4 "4.070.22"
Set-TextValue 16 5 "  +5.14%  "
Set-TextValue 17 4 "21.34"
Set-TextValue 17 5 "  +4.93%  "
Set-TextValue 18 4 "14.39"
Set-TextValue 18 5 "  +3.78%  "
Set-TextValue 19 5 "  +1.63%  "
Set-TextValue 21 4 "72.597.55"
Set-TextValue 21 5 "  +5.58%  "
Set-TextValue 22 4 "455.24"
Set-TextValue 22 5 "  +8.33%  "
Set-TextValue 23 4 "105.24"
Set-TextValue 23 5 "  +21.16%  "
Set-TextValue 24 4 "3.62"
Set-TextValue 24 5 "  +6.99%  "
Set-TextValue 25 4 "15.16"
Set-TextValue 25 5 "  +7.84%  "
Set-TextValue 26 5 "  +1.90%  "
Set-TextValue 27 4 "11.41"
Set-TextValue 27 5 "  +0.83%  "
Set-TextValue 28 5 "  +5.94%  "
Set-TextValue 29 4 "38.27"
Set-TextValue 29 5 "  +6.53%  "
Set-TextValue 30 4 "5.84"
Set-TextValue 30 5 "  +3.17%  "
Set-TextValue 31 5 "  +16.80%  "
Set-TextValue 32 4 "13.70"
Set-TextValue 32 5 "  +4.88%  "
Set-TextValue 33 5 "  +4.66%  "
Set-TextValue 34 4 "677.84"
Set-TextValue 34 5 "  -1.02%  "
Set-TextValue 35 2 "OKB"
Set-TextValue 35 3 "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue 35 4 "67.74"
Set-TextValue 35 5 "  +0.32%  "
Set-TextValue 36 2 "NEARProtocol"
Set-TextValue 36 3 "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue 36 4 "6.63"
Set-TextValue 36 5 "  +12.63%  "
Set-TextValue 37 5 "  +7.09%  "
Set-TextValue 38 2 "PEPE"
Set-TextValue 38 3 "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue 38 4 "0.0₃0868"
Set-TextValue 38 5 "  +2.32%  "
Set-TextValue 39 2 "TheGraph"
Set-TextValue 39 3 "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue 39 4 "0.432"
Set-TextValue 39 5 "  +0.28%  "
Set-TextValue 40 5 "  +4.36%  "
Set-TextValue 41 4 "3.47"
Set-TextValue 41 5 "  +7.64%  "
Set-TextValue 42 4 "1.00"
Set-TextValue 42 5 "  +0.11%  "
Set-TextValue 43 4 "0.0500"
Set-TextValue 43 5 "  +4.87%  "
Set-TextValue 44 5 "  -0.21%  "
Set-TextValue 45 2 "WEMIXToken"
Set-TextValue 45 3 "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue 45 4 "3.22"
Set-TextValue 45 5 "  +2.42%  "
Set-TextValue 46 2 "Stellar"
Set-TextValue 46 3 "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue 46 4 "0.158"
Set-TextValue 46 5 "  +13.44%  "
Set-TextValue 47 4 "2.69"
Set-TextValue 47 5 "  -2.20%  "
Set-TextValue 48 2 "ApeXProtocol"
Set-TextValue 48 3 "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextValue 48 4 "3.48"
Set-TextValue 48 5 "  +2.97%  "
Set-TextValue 49 2 "THORChain"
Set-TextValue 49 3 "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextValue 49 4 "9.81"
Set-TextValue 49 5 "  +16.63%  "
Set-TextValue 50 4 "3.09"
Set-TextValue 50 5 "  +4.78%  "
Set-TextValue 51 2 "LidoDAOToken"
Set-TextValue 51 3 "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue 51 4 "3.36"
Set-TextValue 51 5 "  +3.48%  "
